$wb = $excel.ActiveWorkbook

$wsCustomer   = $wb.Worksheets.Item("Customer")
$wsSuperAdmin = $wb.Worksheets.Item("SuperAdmin")

# --- Customer sheet: move selection/active cell to A11 (no value changes) ---
$wsCustomer.Activate()
$wsCustomer.Range("A11").Select()

# --- SuperAdmin sheet: pick a different dropdown value for A2, select A2,
#     and make SuperAdmin the active (tab-selected) sheet, which also moves
#     the "Repairer" sheet out of being tab-selected. ---
$wsSuperAdmin.Activate()
$wsSuperAdmin.Range("A2").Value = "superadmnphase1@owleyes.ch"
$wsSuperAdmin.Range("A2").Select()
